# Insert a new "Match ID" column at the very left of the sheet (column A),
# shifting all existing columns one place to the right (A:W -> B:X).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()

# Row 2 holds the column headers -> label the new column "Match ID".
$hdr = $ws.Range("A2")
$hdr.Value = "Match ID"
$hdr.Font.Bold = $true

# Data rows (4 through 18) all belong to match id 2.
$data = $ws.Range("A4:A18")
$data.Value = 2
$data.Font.Bold = $true

# The blank spacer / subtotal rows (3, 19) and the new trailing blank row
# (20) just pick up the same bold styling with no value, matching the
# rest of column A.
$ws.Range("A3").Font.Bold = $true
$ws.Range("A19").Font.Bold = $true
$ws.Range("A20").Font.Bold = $true

# Restore the view: keep row 2 pinned at the top, move the active
# selection to E27.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("E27").Select() | Out-Null
